$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 833.3333
$ws.Range("I19").Value = 650
$ws.Range("J19").Value = 1200
$ws.Range("K19").Value = 650
$ws.Range("L19").Value = 1200
$ws.Range("M19").Value = -475
$ws.Range("N19").Value = -1550
$ws.Range("H58").Value = 1149.7368
$ws.Range("I58").Value = 563.4286
$ws.Range("K58").Value = 1690.2858
$ws.Range("M58").Value = -1540.2858
$ws.Range("H98").Value = 1043.8
$ws.Range("I98").Value = 648.1667
$ws.Range("K98").Value = 648.1667
$ws.Range("M98").Value = 849.8333
$ws.Range("H100").Value = 1425.25
$ws.Range("I100").Value = 1400.5
$ws.Range("J100").Value = 1499.5
$ws.Range("K100").Value = 1400.5
$ws.Range("L100").Value = 1499.5
$ws.Range("M100").Value = -859.5
$ws.Range("N100").Value = -2581.5
$ws.Range("H115").Value = 233.75
$ws.Range("I115").Value = 233.75
$ws.Range("K115").Value = 701.25
$ws.Range("M115").Value = 865.75
$ws.Range("H122").Value = 1043.8
$ws.Range("I122").Value = 648.1667
$ws.Range("K122").Value = 1944.5001
$ws.Range("M122").Value = 505.4999
$ws.Range("H132").Value = 3172.7896
$ws.Range("I132").Value = 1138.3636
$ws.Range("J132").Value = 16600
$ws.Range("K132").Value = 3415.0908
$ws.Range("L132").Value = 49800
$ws.Range("M132").Value = -885.0907999999999
$ws.Range("N132").Value = -54860

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 885.3913
$ws.Range("I32").Value = 795.7
$ws.Range("K32").Value = 795.7
$ws.Range("M32").Value = -508.7
$ws.Range("H45").Value = 3799.182
$ws.Range("J45").Value = 4390.625
$ws.Range("L45").Value = 4390.625
$ws.Range("N45").Value = -5144.625
$ws.Range("H74").Value = 2992.6155
$ws.Range("I74").Value = 2298.182
$ws.Range("K74").Value = 2298.182
$ws.Range("M74").Value = -1424.182
$ws.Range("H77").Value = 2992.6155
$ws.Range("I77").Value = 2298.182
$ws.Range("K77").Value = 11490.91
$ws.Range("M77").Value = -7122.91
$ws.Range("H97").Value = 1242.4117
$ws.Range("I97").Value = 1208.4667
$ws.Range("J97").Value = 1497
$ws.Range("K97").Value = 1208.4667
$ws.Range("L97").Value = 1497
$ws.Range("M97").Value = -712.4666999999999
$ws.Range("N97").Value = -2489
$ws.Range("H110").Value = 1763.625
$ws.Range("I110").Value = 1585.1428
$ws.Range("K110").Value = 1585.1428
$ws.Range("M110").Value = 459.8571999999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1222
$ws.Range("I94").Value = 1049.875
$ws.Range("K94").Value = 1049.875
$ws.Range("M94").Value = -598.875

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 104.625
$ws.Range("I7").Value = 58.5
$ws.Range("J7").Value = 150.75
$ws.Range("K7").Value = 58.5
$ws.Range("L7").Value = 150.75
$ws.Range("M7").Value = 54.5
$ws.Range("N7").Value = -376.75
$ws.Range("H31").Value = 4622.385
$ws.Range("I31").Value = 1261.375
$ws.Range("K31").Value = 1261.375
$ws.Range("M31").Value = -966.375
$ws.Range("H34").Value = 4622.385
$ws.Range("I34").Value = 1261.375
$ws.Range("K34").Value = 1261.375
$ws.Range("M34").Value = -1059.375
$ws.Range("H88").Value = 17499.5
$ws.Range("I88").Value = 0
$ws.Range("J88").Value = 17499.5
$ws.Range("K88").Value = 0
$ws.Range("L88").Value = 17499.5
$ws.Range("M88").ClearContents()
$ws.Range("N88").Value = -18311.5
$ws.Range("H91").Value = 17499.5
$ws.Range("I91").Value = 0
$ws.Range("J91").Value = 17499.5
$ws.Range("K91").Value = 0
$ws.Range("L91").Value = 17499.5
$ws.Range("M91").ClearContents()
$ws.Range("N91").Value = -20307.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H13").Value = 456.57144
$ws.Range("I13").Value = 500.5
$ws.Range("J13").Value = 439
$ws.Range("K13").Value = 1501.5
$ws.Range("L13").Value = 1317
$ws.Range("M13").Value = -1333.5
$ws.Range("N13").Value = -1653
$ws.Range("H46").Value = 5559.6
$ws.Range("I46").Value = 932.6667
$ws.Range("K46").Value = 2798.0001
$ws.Range("M46").Value = -2707.0001
$ws.Range("H50").Value = 716.6667
$ws.Range("I50").Value = 850
$ws.Range("J50").Value = 450
$ws.Range("K50").Value = 2550
$ws.Range("L50").Value = 1350
$ws.Range("M50").Value = -2069
$ws.Range("N50").Value = -2312
$ws.Range("H53").Value = 716.6667
$ws.Range("I53").Value = 850
$ws.Range("J53").Value = 450
$ws.Range("K53").Value = 2550
$ws.Range("L53").Value = 1350
$ws.Range("M53").Value = -2069
$ws.Range("N53").Value = -2312
$ws.Range("H75").Value = 3086
$ws.Range("I75").Value = 266.5
$ws.Range("K75").Value = 799.5
$ws.Range("M75").Value = 198.5
$ws.Range("H78").Value = 3086
$ws.Range("I78").Value = 266.5
$ws.Range("K78").Value = 2398.5
$ws.Range("M78").Value = 2593.5
$ws.Range("H114").Value = 1184.8
$ws.Range("I114").Value = 1316.4
$ws.Range("K114").Value = 3949.2
$ws.Range("M114").Value = -695.2000000000003
$ws.Range("H131").Value = 1837.8889
$ws.Range("I131").Value = 1257.6666
$ws.Range("J131").Value = 2998.3333
$ws.Range("K131").Value = 3772.9998
$ws.Range("L131").Value = 8994.999899999999
$ws.Range("M131").Value = 1267.0002
$ws.Range("N131").Value = -19074.9999

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 998.875
$ws.Range("J80").Value = 998.6667
$ws.Range("L80").Value = 998.6667
$ws.Range("N80").Value = -2994.6667
$ws.Range("H83").Value = 998.875
$ws.Range("J83").Value = 998.6667
$ws.Range("L83").Value = 4993.3335
$ws.Range("N83").Value = -14977.3335
$ws.Range("H102").Value = 3244.5
$ws.Range("J102").Value = 5050
$ws.Range("L102").Value = 5050
$ws.Range("N102").Value = -8294
$ws.Range("H113").Value = 2500
$ws.Range("I113").Value = 0
$ws.Range("K113").Value = 0
$ws.Range("M113").ClearContents()
$ws.Range("H122").Value = 3283.5833
$ws.Range("I122").Value = 3145.9
$ws.Range("J122").Value = 3972
$ws.Range("K122").Value = 9437.700000000001
$ws.Range("L122").Value = 11916
$ws.Range("M122").Value = -6987.700000000001
$ws.Range("N122").Value = -16816

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1766.6666
$ws.Range("I22").Value = 2100
$ws.Range("J22").Value = 1700
$ws.Range("K22").Value = 2100
$ws.Range("L22").Value = 1700
$ws.Range("M22").Value = -1805
$ws.Range("N22").Value = -2290
$ws.Range("H27").Value = 1766.6666
$ws.Range("I27").Value = 2100
$ws.Range("J27").Value = 1700
$ws.Range("K27").Value = 2100
$ws.Range("L27").Value = 1700
$ws.Range("M27").Value = -1993
$ws.Range("N27").Value = -1914
$ws.Range("H61").Value = 5327.857
$ws.Range("I61").Value = 4132
$ws.Range("J61").Value = 6224.75
$ws.Range("K61").Value = 4132
$ws.Range("L61").Value = 6224.75
$ws.Range("M61").Value = -3930
$ws.Range("N61").Value = -6628.75
$ws.Range("H113").Value = 5327.857
$ws.Range("I113").Value = 4132
$ws.Range("J113").Value = 6224.75
$ws.Range("K113").Value = 4132
$ws.Range("L113").Value = 6224.75
$ws.Range("M113").Value = -1962
$ws.Range("N113").Value = -10564.75
$ws.Range("H122").Value = 3452
$ws.Range("I122").Value = 3404
$ws.Range("K122").Value = 10212
$ws.Range("M122").Value = -7762

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H80").Value = 60624.5
$ws.Range("J80").Value = 60624.5
$ws.Range("L80").Value = 60624.5
$ws.Range("N80").Value = -62620.5
$ws.Range("H83").Value = 60624.5
$ws.Range("J83").Value = 60624.5
$ws.Range("L83").Value = 181873.5
$ws.Range("N83").Value = -191857.5
$ws.Range("H122").Value = 1466.6666
$ws.Range("J122").Value = 1511.8889
$ws.Range("L122").Value = 4535.6667
$ws.Range("N122").Value = -9435.6667
